# Modularized checks, created Severity check column and changed SimulateClick
# and SimulateType checks to also look for SendWindowMessages.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Workflow"
$ws2 = $wb.Worksheets.Item(2)   # "Project"

# ---------------------------------------------------------------------------
# 1) Insert a new "Severity" column (E) on both sheets, shifting the old
#    Explanation/Suggestion columns one slot to the right.
# ---------------------------------------------------------------------------
$ws1.Range("E1").EntireColumn.Insert()
$ws2.Range("E1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new Severity column + header on sheet "Workflow".
# ---------------------------------------------------------------------------
$ws1.Range("E1").Value = "Severity"
$ws1.Range("E2").Value = "Low"
$ws1.Range("E3").Value = "Low"
$ws1.Range("E4").Value = "Medium"
$ws1.Range("E5").Value = "Medium"
$ws1.Range("E6").Value = "Low"
$ws1.Range("E7").Value = "Low"
$ws1.Range("E8").Value = "Medium"
$ws1.Range("E9").Value = "Medium"
$ws1.Range("E10").Value = "High"
$ws1.Range("E11").Value = "Low"
$ws1.Range("E12").Value = "Low"

# ---------------------------------------------------------------------------
# 3) Populate the new Severity column + header on sheet "Project".
# ---------------------------------------------------------------------------
$ws2.Range("E1").Value = "Severity"
$ws2.Range("E2").Value = "Low"

# ---------------------------------------------------------------------------
# 4) Rename the "SimulateClick"/"SimulateType" checks so they also mention
#    SendWindowMessages (rows 8 and 9 on "Workflow").
# ---------------------------------------------------------------------------
$ws1.Range("B8").Value = "Undocumented default click"
$ws1.Range("C8").Value = "Checks\UndocumentedDefaultClick.xaml"
$ws1.Range("F8").Value = "Since they do not depend on the mouse driver, the properties SimulateClick and SendWindowMessages provide a faster and more robust way to perform clicks, so they should be used whenever possible. Alternatively, add an annotation in case the control does not support such properties. For more about input methods, refer to https://studio.uipath.com/docs/ui-automation#section-input-methods"
$ws1.Range("G8").Value = "Use SimulateClick or SendWindowMessages if the target control supports it."

$ws1.Range("B9").Value = "Undocumented default type"
$ws1.Range("C9").Value = "Checks\UndocumentedDefaultType.xaml"
$ws1.Range("F9").Value = "Since they do not depend on the keyboard driver, the properties SimulateType and SendWindowMessages provide a faster and more robust way to perform typing actions, so they should be used whenever possible. Alternatively, add an annotation in case the control does not support such properties. For more about input methods, refer to https://studio.uipath.com/docs/ui-automation#section-input-methods"
$ws1.Range("G9").Value = "Use SimulateType or SendWindowMessages if the target control supports it."

# ---------------------------------------------------------------------------
# 5) Rename "Unreachable activities" -> "Undocumented unreachable activities"
#    (row 12 on "Workflow").
# ---------------------------------------------------------------------------
$ws1.Range("B12").Value = "Undocumented unreachable activities"
$ws1.Range("C12").Value = "Checks\UndocumentedUnreachableActivities.xaml"

# ---------------------------------------------------------------------------
# 6) Row heights for rows 8 & 9 grow to fit the longer explanation text.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(8).RowHeight = 150
$ws1.Rows.Item(9).RowHeight = 150

# ---------------------------------------------------------------------------
# 7) Column widths: D & E share the old "Argument" width, F keeps the old
#    "Explanation" width, G keeps the old "Suggestion" width.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = $ws1.Columns.Item(4).ColumnWidth
$ws1.Columns.Item(6).ColumnWidth = 49.875
$ws1.Columns.Item(7).ColumnWidth = 39.625

$ws2.Columns.Item(5).ColumnWidth = $ws2.Columns.Item(4).ColumnWidth
$ws2.Columns.Item(6).ColumnWidth = 49.875
$ws2.Columns.Item(7).ColumnWidth = 40.125

# ---------------------------------------------------------------------------
# 8) Update view state: scroll position + active selection moved down/right
#    now that a new row of checks/columns exists.
# ---------------------------------------------------------------------------
$ws1.Application.ActiveWindow.ScrollRow = 10
$ws1.Range("E12").Select()

$ws2.Range("E2").Select()

# ---------------------------------------------------------------------------
# 9) Workbook window size grew slightly.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Height = 15840
